# Updated symbol list - apply price/volume/hour refresh to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.99"
$ws.Range("E2").Value = "'6.43%"
$ws.Range("G2").Value = "'6"
$ws.Range("D3").Value = "'32.25"
$ws.Range("E3").Value = "'10.16%"
$ws.Range("G3").Value = "'6"
$ws.Range("D4").Value = "'5.350"
$ws.Range("E4").Value = "'5.22%"
$ws.Range("G4").Value = "'6"
$ws.Range("D5").Value = "'0.07433"
$ws.Range("E5").Value = "'11.31%"
$ws.Range("G5").Value = "'6"
$ws.Range("D6").Value = "'7.781"
$ws.Range("E6").Value = "'5.71%"
$ws.Range("G6").Value = "'6"
$ws.Range("D7").Value = "'3.692"
$ws.Range("E7").Value = "'8.27%"
$ws.Range("G7").Value = "'6"
$ws.Range("D8").Value = "'1.580"
$ws.Range("E8").Value = "'17.02%"
$ws.Range("G8").Value = "'6"
$ws.Range("D9").Value = "'0.9129"
$ws.Range("E9").Value = "'-0.46%"
$ws.Range("G9").Value = "'6"
$ws.Range("D10").Value = "'0.01717"
$ws.Range("E10").Value = "'2,549.75%"
$ws.Range("G10").Value = "'6"
$ws.Range("D11").Value = "'0.1671"
$ws.Range("E11").Value = "'5.29%"
$ws.Range("G11").Value = "'6"
$ws.Range("D12").Value = "'0.07695"
$ws.Range("E12").Value = "'14.81%"
$ws.Range("G12").Value = "'6"
$ws.Range("D13").Value = "'0.08056"
$ws.Range("E13").Value = "'5.03%"
$ws.Range("G13").Value = "'6"
$ws.Range("D14").Value = "'0.03028"
$ws.Range("E14").Value = "'3.04%"
$ws.Range("G14").Value = "'6"
$ws.Range("D15").Value = "'0.09848"
$ws.Range("E15").Value = "'9.62%"
$ws.Range("G15").Value = "'6"
$ws.Range("D16").Value = "'0.001548"
$ws.Range("E16").Value = "'-1.55%"
$ws.Range("G16").Value = "'6"
$ws.Range("D17").Value = "'0.04551"
$ws.Range("E17").Value = "'0.87%"
$ws.Range("G17").Value = "'6"
$ws.Range("D18").Value = "'0.006350"
$ws.Range("E18").Value = "'1.45%"
$ws.Range("G18").Value = "'6"
$ws.Range("D19").Value = "'3.480"
$ws.Range("E19").Value = "'0.77%"
$ws.Range("G19").Value = "'6"
$ws.Range("D20").Value = "'2.238"
$ws.Range("E20").Value = "'0.83%"
$ws.Range("G20").Value = "'6"
$ws.Range("E21").Value = "'1.60%"
$ws.Range("G21").Value = "'6"
$ws.Range("E22").Value = "'1.92%"
$ws.Range("G22").Value = "'6"
$ws.Range("D23").Value = "'4.200"
$ws.Range("E23").Value = "'3.20%"
$ws.Range("G23").Value = "'6"
$ws.Range("D24").Value = "'0.1621"
$ws.Range("E24").Value = "'3.37%"
$ws.Range("G24").Value = "'6"
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'1.89%"
$ws.Range("G25").Value = "'6"
$ws.Range("D26").Value = "'0.004500"
$ws.Range("E26").Value = "'9.08%"
$ws.Range("G26").Value = "'6"
$ws.Range("D27").Value = "'0.0001170"
$ws.Range("E27").Value = "'-6.38%"
$ws.Range("G27").Value = "'6"
$ws.Range("D28").Value = "'0.0001741"
$ws.Range("E28").Value = "'7.54%"
$ws.Range("G28").Value = "'6"
$ws.Range("G29").Value = "'6"
$ws.Range("G30").Value = "'6"
$ws.Range("G31").Value = "'6"
$ws.Range("G32").Value = "'6"
$ws.Range("G33").Value = "'6"
$ws.Range("G34").Value = "'6"
$ws.Range("G35").Value = "'6"
$ws.Range("G36").Value = "'6"
$ws.Range("G37").Value = "'6"
$ws.Range("G38").Value = "'6"
$ws.Range("G39").Value = "'6"
$ws.Range("D40").Value = "'0.04508"
$ws.Range("E40").Value = "'7.12%"
$ws.Range("G40").Value = "'6"
$ws.Range("D41").Value = "'0.007148"
$ws.Range("E41").Value = "'6.43%"
$ws.Range("G41").Value = "'6"
$ws.Range("D42").Value = "'0.1362"
$ws.Range("E42").Value = "'9.97%"
$ws.Range("G42").Value = "'6"
$ws.Range("D43").Value = "'0.002260"
$ws.Range("E43").Value = "'7.65%"
$ws.Range("G43").Value = "'6"
$ws.Range("D44").Value = "'0.01365"
$ws.Range("E44").Value = "'2.38%"
$ws.Range("G44").Value = "'6"
$ws.Range("D45").Value = "'0.00006133"
$ws.Range("E45").Value = "'6.94%"
$ws.Range("G45").Value = "'6"
$ws.Range("D46").Value = "'1.893"
$ws.Range("E46").Value = "'-3.82%"
$ws.Range("G46").Value = "'6"
$ws.Range("D47").Value = "'0.01300"
$ws.Range("E47").Value = "'-0.55%"
$ws.Range("G47").Value = "'6"
$ws.Range("G48").Value = "'6"
$ws.Range("G49").Value = "'6"
$ws.Range("G50").Value = "'6"
$ws.Range("G51").Value = "'6"

# Re-normalize style so the quote-prefix formatting introduced by the
# text-forcing apostrophe above does not stick to the cells (matches
# original plain/default cell styling).
$ws.Range("D2:G51").Style = "Normal"
